# Update the "timestamp" column (Z) values in the Log_Muestras sheet with the
# new run's timestamps, as recorded when the notebook was re-executed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$timestamps = @{
    2  = "2025-10-17T07:09:29.150001"
    3  = "2025-10-17T07:09:29.150001"
    4  = "2025-10-17T07:09:29.150001"
    5  = "2025-10-17T07:09:29.150001"
    6  = "2025-10-17T07:09:29.150001"
    7  = "2025-10-17T07:09:29.150998"
    8  = "2025-10-17T07:09:29.150998"
    9  = "2025-10-17T07:09:29.150998"
    10 = "2025-10-17T07:09:29.150998"
    11 = "2025-10-17T07:09:29.150998"
    12 = "2025-10-17T07:09:29.150998"
    13 = "2025-10-17T07:09:29.151997"
    14 = "2025-10-17T07:09:29.151997"
    15 = "2025-10-17T07:09:29.151997"
    16 = "2025-10-17T07:09:29.216817"
    17 = "2025-10-17T07:09:29.216817"
    18 = "2025-10-17T07:09:29.217816"
    19 = "2025-10-17T07:09:29.217816"
    20 = "2025-10-17T07:09:29.217816"
    21 = "2025-10-17T07:09:29.217816"
    22 = "2025-10-17T07:09:29.217816"
    23 = "2025-10-17T07:09:29.217816"
    24 = "2025-10-17T07:09:29.218816"
    25 = "2025-10-17T07:09:29.218816"
    26 = "2025-10-17T07:09:29.279467"
    27 = "2025-10-17T07:09:29.279467"
    28 = "2025-10-17T07:09:29.279467"
    29 = "2025-10-17T07:09:29.279467"
    30 = "2025-10-17T07:09:29.279467"
    31 = "2025-10-17T07:09:29.279467"
    32 = "2025-10-17T07:09:29.279467"
    33 = "2025-10-17T07:09:29.279467"
    34 = "2025-10-17T07:09:29.279467"
    35 = "2025-10-17T07:09:29.279467"
    36 = "2025-10-17T07:09:29.279467"
    37 = "2025-10-17T07:09:29.279467"
    38 = "2025-10-17T07:09:29.279467"
    39 = "2025-10-17T07:09:29.279467"
    40 = "2025-10-17T07:09:29.279467"
    41 = "2025-10-17T07:09:29.279467"
    42 = "2025-10-17T07:09:29.279467"
    43 = "2025-10-17T07:09:29.279467"
    44 = "2025-10-17T07:09:29.279467"
    45 = "2025-10-17T07:09:29.279467"
    46 = "2025-10-17T07:09:29.279467"
    47 = "2025-10-17T07:09:29.279467"
    48 = "2025-10-17T07:09:29.279467"
}

foreach ($row in $timestamps.Keys) {
    $ws.Range("Z$row").Value = $timestamps[$row]
}
